$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, pushing the existing row 28 (and below) down to row 29.
$ws.Rows.Item(28).Insert()

# Fill in the new row 28 with the new weekly price record.
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44628
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112040
$ws.Range("G28").Value = "Cilantro"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 240
$ws.Range("K28").Value = 550
$ws.Range("L28").Value = 600
$ws.Range("M28").Value = 575
$ws.Range("N28").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O28").Value = "Provincia de Diguillín"
$ws.Range("P28").Value = 575
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"

# Ensure D28 keeps the same date style (s="2") as the rest of the date column.
$ws.Range("D28").NumberFormat = $ws.Range("D27").NumberFormat
